# Update column F (dSF) values to reflect repulled data / recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 7
$ws.Range("F6").Value = 1
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -4
$ws.Range("F14").Value = -5
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = -4
